# Add a new "采购价" (purchase price) column to the goods sheet, inserted
# right before the existing "库存" (stock) column, and populate it with
# values for each product row. Also restore the active selection to B8
# (as left by the author after the edit).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at G, shifting 库存/状态/详情 one column to the right.
$ws.Columns("G:G").Insert()

# Header for the new column.
$ws.Range("G1").Value = "采购价"

# Purchase-price values for each of the four product rows.
$ws.Range("G2").Value = 88
$ws.Range("G3").Value = 6.8
$ws.Range("G4").Value = 0.01
$ws.Range("G5").Value = 0.01

# Restore selection left by the author.
$ws.Range("B8").Select()
